$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 69623867
$ws.Range("B2").Value = 77506
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 400592.8898661289
$ws.Range("R2").Value = 6739831.101304529

# Row 3
$ws.Range("A3").Value = 69623851
$ws.Range("B3").Value = 89392
$ws.Range("E3").Value = 1202
$ws.Range("F3").Value = "Ullticka"
$ws.Range("G3").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H3").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q3").Value = 400558.1253763044
$ws.Range("R3").Value = 6739813.010787125

# Row 4
$ws.Range("A4").Value = 69623843
$ws.Range("B4").Value = 78098
$ws.Range("E4").Value = 6453
$ws.Range("F4").Value = "Vedskivlav"
$ws.Range("G4").Value = "Hertelidea botryosa"
$ws.Range("H4").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q4").Value = 400160.2169516565
$ws.Range("R4").Value = 6739725.882952046

# Row 5
$ws.Range("A5").Value = 69623846
$ws.Range("B5").Value = 77506
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 400325.1187917067
$ws.Range("R5").Value = 6739724.211813147

# Row 6
$ws.Range("A6").Value = 69623848
$ws.Range("B6").Value = 77258
$ws.Range("E6").Value = 6446
$ws.Range("F6").Value = "Kolflarnlav"
$ws.Range("G6").Value = "Carbonicola anthracophila"
$ws.Range("H6").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q6").Value = 400325.1187917067
$ws.Range("R6").Value = 6739724.211813147

# Row 7
$ws.Range("A7").Value = 69623847
$ws.Range("B7").Value = 78098
$ws.Range("E7").Value = 6453
$ws.Range("F7").Value = "Vedskivlav"
$ws.Range("G7").Value = "Hertelidea botryosa"
$ws.Range("H7").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q7").Value = 400325.1187917067
$ws.Range("R7").Value = 6739724.211813147
